# Auto-generated edit script applying cryptos.xlsx price/volume update
# (commit: "Updated cryptos list on Tue Aug  6 21:31:30 UTC 2024 with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.184.46"
$ws.Range("E2").Value = "  +3.22%  "
$ws.Range("D3").Value = "2.476.34"
$ws.Range("E3").Value = "  +2.31%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'488.68"
$ws.Range("E5").Value = "  +5.13%  "
$ws.Range("D6").Value = "'147.02"
$ws.Range("E6").Value = "  +10.41%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").Value = "'0.509"
$ws.Range("E8").Value = "  +3.34%  "
$ws.Range("D9").Value = "2.483.36"
$ws.Range("E9").Value = "  +1.80%  "
$ws.Range("D10").Value = "'5.80"
$ws.Range("E10").Value = "  +9.12%  "
$ws.Range("D11").Value = "'0.0970"
$ws.Range("E11").Value = "  +1.97%  "
$ws.Range("D12").Value = "'0.332"
$ws.Range("E12").Value = "  +5.42%  "
$ws.Range("E13").Value = "  +1.37%  "
$ws.Range("D14").Value = "2.908.39"
$ws.Range("E14").Value = "  +1.99%  "
$ws.Range("D15").Value = "56.203.26"
$ws.Range("E15").Value = "  +3.34%  "
$ws.Range("D16").Value = "'21.08"
$ws.Range("E16").Value = "  +6.60%  "
$ws.Range("E17").Value = "  +2.26%  "
$ws.Range("D18").Value = "2.486.09"
$ws.Range("E18").Value = "  +1.51%  "
$ws.Range("D19").Value = "'4.52"
$ws.Range("E19").Value = "  +7.81%  "
$ws.Range("D20").Value = "'10.06"
$ws.Range("E20").Value = "  +6.41%  "
$ws.Range("D21").Value = "'318.18"
$ws.Range("E21").Value = "  +2.67%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  +0.38%  "
$ws.Range("D23").Value = "'5.80"
$ws.Range("E23").Value = "  +7.91%  "
$ws.Range("D24").Value = "'58.37"
$ws.Range("E24").Value = "  +3.79%  "
$ws.Range("E25").Value = "  +7.28%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.52%  "
$ws.Range("E27").Value = "  +4.33%  "
$ws.Range("D28").Value = "2.580.62"
$ws.Range("E28").Value = "  +2.20%  "
$ws.Range("E29").Value = "  +7.04%  "
$ws.Range("D30").Value = "0.0₃0791"
$ws.Range("E30").Value = "  +10.26%  "
$ws.Range("D31").Value = "'0.998"
$ws.Range("E31").Value = "  +0.09%  "
$ws.Range("D32").Value = "'148.95"
$ws.Range("E32").Value = "  +1.27%  "
$ws.Range("D33").Value = "'18.19"
$ws.Range("E33").Value = "  +2.44%  "
$ws.Range("E34").Value = "  +5.11%  "
$ws.Range("D35").Value = "'5.19"
$ws.Range("E35").Value = "  +4.17%  "
$ws.Range("E36").Value = "  +8.19%  "
$ws.Range("D37").Value = "'3.74"
$ws.Range("E37").Value = "  +5.30%  "
$ws.Range("D38").Value = "'0.862"
$ws.Range("E38").Value = "  +6.63%  "
$ws.Range("D39").Value = "'34.20"
$ws.Range("E39").Value = "  +4.10%  "
$ws.Range("D40").Value = "'3.52"
$ws.Range("E40").Value = "  +8.36%  "
$ws.Range("E41").Value = "  +0.12%  "
$ws.Range("D42").Value = "'0.0556"
$ws.Range("E42").Value = "  +5.98%  "
$ws.Range("D43").Value = "'0.604"
$ws.Range("E43").Value = "  +1.40%  "
$ws.Range("D44").Value = "'1.33"
$ws.Range("E44").Value = "  +7.26%  "
$ws.Range("D45").Value = "'4.80"
$ws.Range("E45").Value = "  +13.82%  "
$ws.Range("D46").Value = "'0.0925"
$ws.Range("E46").Value = "  +4.65%  "
$ws.Range("B47").Value = "WhiteBITCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D47").Value = "'10.20"
$ws.Range("E47").Value = "  +1.26%  "
$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").Value = "'257.46"
$ws.Range("E48").Value = "  +11.62%  "
$ws.Range("D49").Value = "'0.0228"
$ws.Range("E49").Value = "  +4.79%  "
$ws.Range("D50").Value = "'17.62"
$ws.Range("E50").Value = "  +6.21%  "
$ws.Range("D51").Value = "1.872.52"
$ws.Range("E51").Value = "  -3.54%  "
